$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "<Name>_old" -> "<Name>_FV2404" (columns A:J)
#    and "<Name>_new" -> "<Name>_FV2410" (columns L:U). Column K ("diff") stays the same.
$oldHeaders = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)
$newHeaders = @(
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)

for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# 2. Turn the data range into an Excel Table named "Table1"
$range = $ws.Range("A1:U78")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze the header row (split below row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
